$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.324075666666667
$ws.Range("H2").Value = 3.972227
$ws.Range("I2").Value = 0.01675578032580584
$ws.Range("J2").Value = 0.01684165790066494
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.874784666666667
$ws.Range("N2").Value = 8.624354
$ws.Range("O2").Value = 0.1187109652550681
$ws.Range("P2").Value = 0.121184727686443
$ws.Range("Q2").Value = 3.806432424039778
$ws.Range("R2").Value = 34.257891816358
$ws.Range("S2").Value = 0.00198909485607829
$ws.Range("T2").Value = 0.002040951726480312

$ws.Range("G3").Value = 1.324075666666667
$ws.Range("H3").Value = 3.972227
$ws.Range("I3").Value = 0.01675578032580584
$ws.Range("J3").Value = 0.01684165790066494
$ws.Range("O3").Value = 0.4442422727481699
$ws.Range("P3").Value = 0.4534996302499962
$ws.Range("Q3").Value = 14.24449870729666
$ws.Range("R3").Value = 128.20048836567
$ws.Range("S3").Value = 0.007443625933605056
$ws.Range("T3").Value = 0.007637685630748476

$ws.Range("G4").Value = 1.324075666666667
$ws.Range("H4").Value = 3.972227
$ws.Range("I4").Value = 0.01675578032580584
$ws.Range("J4").Value = 0.01684165790066494
$ws.Range("M4").Value = 5.147441999999999
$ws.Range("N4").Value = 15.442326
$ws.Range("O4").Value = 0.2125577666737049
$ws.Range("P4").Value = 0.2169871588243338
$ws.Range("Q4").Value = 6.815602697777998
$ws.Range("R4").Value = 61.340424280002
$ws.Range("S4").Value = 0.003561571244928493
$ws.Range("T4").Value = 0.003654423497756679

$ws.Range("G5").Value = 1.324075666666667
$ws.Range("H5").Value = 3.972227
$ws.Range("I5").Value = 0.01675578032580584
$ws.Range("J5").Value = 0.01684165790066494
$ws.Range("M5").Value = 1.483016
$ws.Range("N5").Value = 2.966032
$ws.Range("O5").Value = 0.06123946008548931
$ws.Range("P5").Value = 0.04167706708575228
$ws.Range("Q5").Value = 1.963625398877334
$ws.Range("R5").Value = 11.781752393264
$ws.Range("S5").Value = 0.001026114940463414
$ws.Range("T5").Value = 0.0007019109061613025

$ws.Range("G6").Value = 1.324075666666667
$ws.Range("H6").Value = 3.972227
$ws.Range("I6").Value = 0.01675578032580584
$ws.Range("J6").Value = 0.01684165790066494
$ws.Range("M6").Value = 3.953360666666667
$ws.Range("N6").Value = 11.860082
$ws.Range("O6").Value = 0.1632495352375677
$ws.Range("P6").Value = 0.1666514161534747
$ws.Range("Q6").Value = 5.234548660290445
$ws.Range("R6").Value = 47.110937942614
$ws.Range("S6").Value = 0.002735373350730584
$ws.Range("T6").Value = 0.002806686139518168

$ws.Range("I7").Value = 0.9679468703219594
$ws.Range("J7").Value = 0.9729078406975189
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.874784666666667
$ws.Range("N7").Value = 8.624354
$ws.Range("O7").Value = 0.1187109652550681
$ws.Range("P7").Value = 0.121184727686443
$ws.Range("Q7").Value = 219.8897503010882
$ws.Range("R7").Value = 1979.007752709794
$ws.Range("S7").Value = 0.114905907291542
$ws.Range("T7").Value = 0.1179015717389341

$ws.Range("I8").Value = 0.9679468703219594
$ws.Range("J8").Value = 0.9729078406975189
$ws.Range("O8").Value = 0.4442422727481699
$ws.Range("P8").Value = 0.4534996302499962
$ws.Range("R8").Value = 7405.877797059808
$ws.Range("S8").Value = 0.4300029175713053
$ws.Range("T8").Value = 0.441213346023647

$ws.Range("I9").Value = 0.9679468703219594
$ws.Range("J9").Value = 0.9729078406975189
$ws.Range("M9").Value = 5.147441999999999
$ws.Range("N9").Value = 15.442326
$ws.Range("O9").Value = 0.2125577666737049
$ws.Range("P9").Value = 0.2169871588243338
$ws.Range("Q9").Value = 393.7233105468539
$ws.Range("R9").Value = 3543.509794921685
$ws.Range("S9").Value = 0.205744625014438
$ws.Range("T9").Value = 0.2111085081508722

$ws.Range("I10").Value = 0.9679468703219594
$ws.Range("J10").Value = 0.9729078406975189
$ws.Range("M10").Value = 1.483016
$ws.Range("N10").Value = 2.966032
$ws.Range("O10").Value = 0.06123946008548931
$ws.Range("P10").Value = 0.04167706708575228
$ws.Range("Q10").Value = 113.4345892802587
$ws.Range("R10").Value = 680.607535681552
$ws.Range("S10").Value = 0.05927654372995593
$ws.Range("T10").Value = 0.04054794534500489

$ws.Range("I11").Value = 0.9679468703219594
$ws.Range("J11").Value = 0.9729078406975189
$ws.Range("M11").Value = 3.953360666666667
$ws.Range("N11").Value = 11.860082
$ws.Range("O11").Value = 0.1632495352375677
$ws.Range("P11").Value = 0.1666514161534747
$ws.Range("Q11").Value = 302.3890797537335
$ws.Range("R11").Value = 2721.501717783602
$ws.Range("S11").Value = 0.1580168767147181
$ws.Range("T11").Value = 0.1621364694390607

$ws.Range("G12").Value = 1.2088275
$ws.Range("H12").Value = 2.417655
$ws.Range("I12").Value = 0.0152973493522347
$ws.Range("J12").Value = 0.01025050140181618
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.874784666666667
$ws.Range("N12").Value = 8.624354
$ws.Range("O12").Value = 0.1187109652550681
$ws.Range("P12").Value = 0.121184727686443
$ws.Range("Q12").Value = 3.475118761645
$ws.Range("R12").Value = 20.85071256987
$ws.Range("S12").Value = 0.001815963107447771
$ws.Range("T12").Value = 0.001242204221028597

$ws.Range("G13").Value = 1.2088275
$ws.Range("H13").Value = 2.417655
$ws.Range("I13").Value = 0.0152973493522347
$ws.Range("J13").Value = 0.01025050140181618
$ws.Range("O13").Value = 0.4442422727481699
$ws.Range("P13").Value = 0.4534996302499962
$ws.Range("Q13").Value = 13.004650862925
$ws.Range("R13").Value = 78.02790517754998
$ws.Range("S13").Value = 0.006795729243259485
$ws.Range("T13").Value = 0.004648598595600706

$ws.Range("G14").Value = 1.2088275
$ws.Range("H14").Value = 2.417655
$ws.Range("I14").Value = 0.0152973493522347
$ws.Range("J14").Value = 0.01025050140181618
$ws.Range("M14").Value = 5.147441999999999
$ws.Range("N14").Value = 15.442326
$ws.Range("O14").Value = 0.2125577666737049
$ws.Range("P14").Value = 0.2169871588243338
$ws.Range("Q14").Value = 6.222369444254999
$ws.Range("R14").Value = 37.33421666552999
$ws.Range("S14").Value = 0.003251570414338454
$ws.Range("T14").Value = 0.002224227175704945

$ws.Range("G15").Value = 1.2088275
$ws.Range("H15").Value = 2.417655
$ws.Range("I15").Value = 0.0152973493522347
$ws.Range("J15").Value = 0.01025050140181618
$ws.Range("M15").Value = 1.483016
$ws.Range("N15").Value = 2.966032
$ws.Range("O15").Value = 0.06123946008548931
$ws.Range("P15").Value = 0.04167706708575228
$ws.Range("Q15").Value = 1.79271052374
$ws.Range("R15").Value = 7.17084209496
$ws.Range("S15").Value = 0.0009368014150699625
$ws.Range("T15").Value = 0.0004272108345860909

$ws.Range("G16").Value = 1.2088275
$ws.Range("H16").Value = 2.417655
$ws.Range("I16").Value = 0.0152973493522347
$ws.Range("J16").Value = 0.01025050140181618
$ws.Range("M16").Value = 3.953360666666667
$ws.Range("N16").Value = 11.860082
$ws.Range("O16").Value = 0.1632495352375677
$ws.Range("P16").Value = 0.1666514161534747
$ws.Range("Q16").Value = 4.778931091285
$ws.Range("R16").Value = 28.67358654771
$ws.Range("S16").Value = 0.002497285172119022
$ws.Range("T16").Value = 0.001708260574895845
